$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells P1 and Q1, copying the formatting (bold, border, centered) from O1
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Update the full data block B2:Q25 with the new dataset values (new columns P and Q added)
$data = @(
    @(19.51318330836458, 0, 7.625465523357247, 13.58676995102795, 31.64194912159061, 48.54607971675767, 1.609263798000339, 3.586837643496082, 13.63679600169964, 15.73501389633077, 8.671892546297096, 0, 13.96630972077187, 0, 13.48489331570636, 0),
    @(18.33313029076075, 0, 7.297316953275867, 12.96573967958801, 30.31210435103559, 45.79584229638348, 1.914737674834541, 3.743678253478187, 13.16381228713522, 15.93135596563163, 8.378348064866142, 0, 13.13127815516542, 0, 13.51130099949605, 0),
    @(17.5662943636996, 0, 7.090160009296527, 12.57233218548465, 29.48332758302094, 44.05225463426169, 2.107785136851385, 3.843411090631794, 12.87264200722417, 16.05536014986594, 8.193494371581171, 0, 12.59552240399426, 0, 13.53038950290296, 0),
    @(17.23515459419759, 0, 7.005115923618462, 12.41221589806561, 29.13372356115362, 43.3149742587939, 2.188139330615337, 3.887512499780473, 12.75103473069179, 16.10131798039138, 8.1185131505099, 0, 12.37710404487215, 0, 13.54066087898611, 0),
    @(17.16960306779956, 0, 6.991878177523153, 12.38939283970173, 29.06427547545983, 43.17517577142967, 2.201888644372392, 3.898185379555805, 12.72730161863119, 16.10238729060502, 8.107834426369244, 0, 12.34734323615923, 0, 13.54460103761869, 0),
    @(17.53533203138073, 0, 7.091636713312051, 12.580821190824, 29.44836002775223, 43.99796929634051, 2.10971541751759, 3.852722756938251, 12.86144953838545, 16.03790179416024, 8.197487407692741, 0, 12.61087653970447, 0, 13.53662037149287, 0),
    @(19.08251437441515, 0, 7.516828199032721, 13.38859428103636, 31.14821611800775, 47.55487019980743, 1.714522547757771, 3.651280181446871, 13.46185082890347, 15.77803394645077, 8.578110400985796, 0, 13.70561189224539, 0, 13.50194570936554, 0),
    @(21.82065879529994, 0, 8.301322424219761, 14.8589180944074, 34.42003487939468, 54.14588104484299, 2.127539109173018, 3.26841209373434, 14.64430064020431, 15.32933303808247, 9.279230495441288, 0, 15.63725151701856, 0, 13.44551265990597, 0),
    @(23.54405280449025, 0, 8.872012353151916, 15.96464837536847, 36.53300192521976, 58.402183840493, 2.621074761318226, 3.014241657287162, 15.43382288007927, 14.95364239468854, 9.8162571738129, 0, 16.94541007019691, 0, 13.44975237240592, 0),
    @(23.65901271111852, 0, 9.300968800903178, 17.16487185448434, 35.75872294686096, 57.87570533050281, 3.402587053058, 2.974660627748363, 15.23303078796699, 14.29807557538748, 10.44255369121973, 0, 16.68711847595544, 0, 13.72402872918187, 0),
    @(23.42955321165781, 0, 9.546191077911155, 17.97114764481506, 34.67614002656583, 56.57667289682337, 4.535405010519355, 2.97498008325116, 14.9004099388142, 13.88800006374782, 10.88195871093604, 0, 16.17478823793232, 0, 13.9648630335003, 0),
    @(22.88328461991764, 0, 9.672197969120591, 18.55065388983255, 33.21773561949743, 54.49653604750739, 5.762194475771373, 3.020870818068714, 14.42325456157468, 13.62159429727174, 11.21442014242361, 0, 15.44050443504515, 0, 14.19849818160294, 0),
    @(22.34359724476528, 0, 9.709120461007766, 18.86634074993401, 32.00263219788223, 52.64620615474004, 6.657499048449204, 3.076289052289299, 14.01677049300663, 13.50701962835525, 11.40538725618028, 0, 14.81999368058788, 0, 14.36143278553851, 0),
    @(22.13469650519854, 0, 9.693313878591399, 18.89641874472796, 31.61619983272191, 52.00989075585773, 6.859231870542684, 3.104125233088976, 13.88405472732831, 13.5019638778287, 11.42853674504503, 0, 14.62416466984283, 0, 14.3999622510133, 0),
    @(21.48459085893023, 0, 9.43756228659395, 18.36203483471355, 30.90845158570686, 50.45339306889171, 6.563368224155903, 3.209256879184845, 13.61292030598704, 13.7185714647449, 11.15273334930718, 0, 14.180555968554, 0, 14.33900349180367, 0),
    @(21.27651685050012, 0, 9.217774786547459, 17.76686362323173, 31.03636872887907, 50.28263929554939, 5.806768952770586, 3.261139577858025, 13.62862202035429, 13.93797543535397, 10.83115690378723, 0, 14.1824318780294, 0, 14.20711260847108, 0),
    @(21.45198882858519, 0, 9.004937784032197, 17.0581339223964, 31.91365084614631, 51.32192885511528, 4.607189858720228, 3.261707724079431, 13.89653219546797, 14.21744345847605, 10.440692592346, 0, 14.56452337239651, 0, 13.999640151091, 0),
    @(21.87793686875108, 0, 8.822989662117662, 16.35979787312849, 33.25920504197526, 53.16094945175696, 3.303102847894861, 3.237076928565912, 14.32471113708573, 14.53733090112833, 10.05603655375137, 0, 15.24565206172459, 0, 13.77558317026433, 0),
    @(23.03696497508014, 0, 8.731895898035912, 15.71026274923874, 35.89638802477525, 57.18177778041245, 2.489097887602304, 3.109849768569878, 15.1991111964469, 14.99231748168816, 9.693011421543723, 0, 16.64737822928773, 0, 13.4713601552941, 0),
    @(24.4182496636974, 0, 9.129738905873634, 16.41285524101758, 37.81165466770868, 60.80929469475171, 2.881715640768495, 2.903540418006242, 15.90421060417931, 14.80905397248131, 10.03074475312031, 0, 17.75483452782389, 0, 13.43347357422233, 0),
    @(25.27180955583727, 0, 9.393062580123852, 16.8923659029111, 38.97519877448515, 63.02474324295253, 3.129559411012325, 2.762119063202865, 16.33923251207548, 14.68923110552413, 10.2637311432176, 0, 18.3861707002807, 0, 13.41566851698713, 0),
    @(24.84283699918909, 0, 9.250359039858603, 16.62719149615314, 38.38370429561216, 61.88555062963269, 2.998224598009378, 2.825558675559368, 16.11645165612621, 14.77446804077141, 10.13450020118874, 0, 18.03624273781609, 0, 13.41606805761553, 0),
    @(23.10698930531235, 0, 8.705543271393021, 15.61966135855807, 36.06618456548972, 57.40481221901292, 2.494186925507773, 3.089734047775005, 15.25259955214036, 15.07224005925181, 9.646009480698368, 0, 16.68385802374906, 0, 13.4373568115589, 0),
    @(21.07600031799193, 0, 8.098117755693458, 14.48852028036628, 33.50060931183064, 52.34609614998803, 1.93627377482939, 3.385148285772308, 14.31074069504835, 15.41470379935869, 9.101970865250184, 0, 15.16255774834802, 0, 13.47071750546989, 0)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowIndex = $i + 2
    for ($j = 0; $j -lt $data[$i].Count; $j++) {
        $colIndex = $j + 2  # column B = 2
        $ws.Cells.Item($rowIndex, $colIndex).Value = $data[$i][$j]
    }
}
